# "ammended the terraform slide and script"
#
# On the Terraform slide, after the existing bullet that ends with
# "...vms in AWS.", add a blank line followed by a new paragraph describing
# the inventory.yaml output:
#
#   Outputs an inventory.yaml containing vm IP addresses.

$p = $ppt.ActivePresentation

# Locate the slide / shape holding the Terraform bullet list rather than
# hard-coding indices, so the script is resilient to slide re-ordering.
$anchor = "Used it to set up and configure our vms in AWS."
$targetSlide = $null
$targetShape = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text -like "*$anchor*") {
                $targetSlide = $slide
                $targetShape = $shape
                break
            }
        }
    }
    if ($targetShape -ne $null) { break }
}

$tr = $targetShape.TextFrame.TextRange

# Calling InsertAfter repeatedly on the *original* whole-body TextRange (not
# on a narrowed sub-range) appends each new chunk as its own run right after
# the previous one while leaving every pre-existing paragraph byte-for-byte
# untouched. A CR ("`r") is PowerPoint's paragraph-separator inside
# TextRange.Text/InsertAfter.

# 1) a new, empty paragraph between the old bullet and the new one
$tr.InsertAfter("`r`r") | Out-Null

# 2) the new paragraph, built up phrase by phrase so each becomes its own run
#    (mirrors how the author's edit was split across multiple <a:r> runs)
$tr.InsertAfter("Outputs an ") | Out-Null
$tr.InsertAfter("inventory.yaml") | Out-Null
$tr.InsertAfter(" containing vm ") | Out-Null
$tr.InsertAfter("IP addresses.") | Out-Null
